$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 6275.2
$ws.Range("I54").Value = 6275.2
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 6275.2
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -5789.2

$ws.Range("H62").Value = 2412.2
$ws.Range("I62").Value = 1814
$ws.Range("J62").Value = 4805
$ws.Range("K62").Value = 1814
$ws.Range("L62").Value = 4805
$ws.Range("M62").Value = -1190
$ws.Range("N62").Value = -6053

$ws.Range("H65").Value = 2412.2
$ws.Range("I65").Value = 1814
$ws.Range("J65").Value = 4805
$ws.Range("K65").Value = 9070
$ws.Range("L65").Value = 24025
$ws.Range("M65").Value = -5950
$ws.Range("N65").Value = -30265

$ws.Range("H86").Value = 2271.1428
$ws.Range("J86").Value = 1500
$ws.Range("L86").Value = 1500
$ws.Range("N86").Value = -3746

$ws.Range("H89").Value = 2271.1428
$ws.Range("J89").Value = 1500
$ws.Range("L89").Value = 7500
$ws.Range("N89").Value = -18732

$ws.Range("H98").Value = 2694.6667
$ws.Range("J98").Value = 3488.5
$ws.Range("L98").Value = 3488.5
$ws.Range("N98").Value = -6484.5

$ws.Range("H116").Value = 7898.5
$ws.Range("I116").Value = 7898.5
$ws.Range("K116").Value = 7898.5
$ws.Range("M116").Value = -4456.5

$ws.Range("H122").Value = 2694.6667
$ws.Range("J122").Value = 3488.5
$ws.Range("L122").Value = 10465.5
$ws.Range("N122").Value = -15365.5

$ws.Range("H132").Value = 3157.963
$ws.Range("I132").Value = 2699.0833
$ws.Range("K132").Value = 8097.249899999999
$ws.Range("M132").Value = -5567.249899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2007.0714
$ws.Range("I2").Value = 2156.2727
$ws.Range("J2").Value = 1460
$ws.Range("K2").Value = 2156.2727
$ws.Range("L2").Value = 1460
$ws.Range("M2").Value = -2043.2727
$ws.Range("N2").Value = -1686

$ws.Range("H32").Value = 3103.3901
$ws.Range("I32").Value = 1769.7894
$ws.Range("K32").Value = 1769.7894
$ws.Range("M32").Value = -1482.7894

$ws.Range("H43").Value = 126666
$ws.Range("J43").Value = 39999.5
$ws.Range("L43").Value = 39999.5
$ws.Range("N43").Value = -40625.5

$ws.Range("H116").Value = 2007.0714
$ws.Range("I116").Value = 2156.2727
$ws.Range("J116").Value = 1460
$ws.Range("K116").Value = 2156.2727
$ws.Range("L116").Value = 1460
$ws.Range("M116").Value = 137.7273
$ws.Range("N116").Value = -6048

$ws.Range("H132").Value = 2061.1177
$ws.Range("I132").Value = 1782.7
$ws.Range("K132").Value = 5348.1
$ws.Range("M132").Value = -2818.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2007.0714
$ws.Range("I3").Value = 2156.2727
$ws.Range("J3").Value = 1460
$ws.Range("K3").Value = 2156.2727
$ws.Range("L3").Value = 1460
$ws.Range("M3").Value = -2042.2727
$ws.Range("N3").Value = -1688

$ws.Range("H26").Value = 13787.8
$ws.Range("I26").Value = 13787.8
$ws.Range("K26").Value = 13787.8
$ws.Range("M26").Value = -13495.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 39792.8
$ws.Range("I59").Value = 27498.5
$ws.Range("J59").Value = 47989
$ws.Range("K59").Value = 27498.5
$ws.Range("L59").Value = 47989
$ws.Range("M59").Value = -26353.5
$ws.Range("N59").Value = -50279

$ws.Range("H60").Value = 20632.334
$ws.Range("I60").Value = 20632.334
$ws.Range("K60").Value = 20632.334
$ws.Range("M60").Value = -20121.334

$ws.Range("H122").Value = 1830
$ws.Range("I122").Value = 1720.8462
$ws.Range("J122").Value = 2184.75
$ws.Range("K122").Value = 5162.5386
$ws.Range("L122").Value = 6554.25
$ws.Range("M122").Value = -2712.5386
$ws.Range("N122").Value = -11454.25

$ws.Range("H132").Value = 4149.9165
$ws.Range("I132").Value = 3980.3
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 11940.9
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -9410.900000000001
$ws.Range("N132").Value = -20054

$ws.Range("H134").Value = 3750
$ws.Range("I134").Value = 2454.6365
$ws.Range("J134").Value = 17999
$ws.Range("K134").Value = 7363.9095
$ws.Range("L134").Value = 53997
$ws.Range("M134").Value = -4828.9095
$ws.Range("N134").Value = -59067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 900
$ws.Range("L22").Value = 2700
$ws.Range("N22").Value = -3038

$ws.Range("H27").Value = 1000
$ws.Range("J27").Value = 900
$ws.Range("L27").Value = 2700
$ws.Range("N27").Value = -2904

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = $null

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = $null

$ws.Range("H68").Value = 3855.375
$ws.Range("I68").Value = 4290.6665
$ws.Range("K68").Value = 12871.9995
$ws.Range("M68").Value = -12060.9995

$ws.Range("H71").Value = 3855.375
$ws.Range("I71").Value = 4290.6665
$ws.Range("K71").Value = 38615.9985
$ws.Range("M71").Value = -34559.9985

$ws.Range("H122").Value = 460
$ws.Range("I122").Value = 460
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4140
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null
$ws.Range("M122").Value = -1690

$ws.Range("H129").Value = 3822.5334
$ws.Range("I129").Value = 996.6667
$ws.Range("J129").Value = 4529
$ws.Range("K129").Value = 2990.0001
$ws.Range("L129").Value = 13587
$ws.Range("M129").Value = 2009.9999
$ws.Range("N129").Value = -23587

$ws.Range("H132").Value = 6165.5
$ws.Range("I132").Value = 4997
$ws.Range("J132").Value = 6399.2
$ws.Range("K132").Value = 44973
$ws.Range("L132").Value = 57592.8
$ws.Range("M132").Value = -42443
$ws.Range("N132").Value = -62652.8

$ws.Range("H136").Value = 5620
$ws.Range("I136").Value = 4028.5715
$ws.Range("K136").Value = 12085.7145
$ws.Range("M136").Value = -6985.7145

$ws.Range("H138").Value = 1713.2
$ws.Range("I138").Value = 1713.2
$ws.Range("K138").Value = 5139.6
$ws.Range("M138").Value = 0.3999999999996362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3204.3684
$ws.Range("I126").Value = 2668.7693
$ws.Range("J126").Value = 4364.8335
$ws.Range("K126").Value = 8006.3079
$ws.Range("L126").Value = 13094.5005
$ws.Range("M126").Value = -5536.3079
$ws.Range("N126").Value = -18034.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4232.8887
$ws.Range("I7").Value = 1699.6666
$ws.Range("K7").Value = 1699.6666
$ws.Range("M7").Value = -1587.6666

$ws.Range("H40").Value = 7102.923
$ws.Range("I40").Value = 7556.5
$ws.Range("J40").Value = 6714.143
$ws.Range("K40").Value = 7556.5
$ws.Range("L40").Value = 6714.143
$ws.Range("M40").Value = -7420.5
$ws.Range("N40").Value = -6986.143

$ws.Range("H46").Value = 1502
$ws.Range("I46").Value = 1002.2857
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 1002.2857
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -814.2857
$ws.Range("N46").Value = -5376

$ws.Range("H93").Value = 3007.8
$ws.Range("I93").Value = 3504.4285
$ws.Range("J93").Value = 1849
$ws.Range("K93").Value = 3504.4285
$ws.Range("L93").Value = 1849
$ws.Range("M93").Value = -2256.4285
$ws.Range("N93").Value = -4345

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null

$ws.Range("H126").Value = 4232.8887
$ws.Range("I126").Value = 1699.6666
$ws.Range("K126").Value = 5098.9998
$ws.Range("M126").Value = -2628.9998

$ws.Range("H132").Value = 3989
$ws.Range("I132").Value = 3613.1428
$ws.Range("J132").Value = 4281.3335
$ws.Range("K132").Value = 10839.4284
$ws.Range("L132").Value = 12844.0005
$ws.Range("M132").Value = -8309.4284
$ws.Range("N132").Value = -17904.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1917.6923
$ws.Range("J81").Value = 3500.3333
$ws.Range("L81").Value = 7000.6666
$ws.Range("N81").Value = -9122.6666

$ws.Range("H84").Value = 1917.6923
$ws.Range("J84").Value = 3500.3333
$ws.Range("L84").Value = 35003.333
$ws.Range("N84").Value = -45611.333

$ws.Range("H99").Value = 53999.5
$ws.Range("I99").Value = 53999.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 53999.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -51004.5
$ws.Range("N99").Value = $null

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null

$ws.Range("H126").Value = 3498.7
$ws.Range("I126").Value = 3564.0625
$ws.Range("J126").Value = 3237.25
$ws.Range("K126").Value = 10692.1875
$ws.Range("L126").Value = 9711.75
$ws.Range("M126").Value = -8222.1875
$ws.Range("N126").Value = -14651.75

$ws.Range("H132").Value = 5200.3
$ws.Range("I132").Value = 5304.074
$ws.Range("J132").Value = 4266.3335
$ws.Range("K132").Value = 15912.222
$ws.Range("L132").Value = 12799.0005
$ws.Range("M132").Value = -13382.222
$ws.Range("N132").Value = -17859.0005
